$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("parameter") is updated per row based on the row's datatype
# (column B): bool -> bool_para, float -> float_para, integer -> int_para,
# string -> str_para, list -> list_para.
#
# The new shared-string entries must be introduced in this exact order
# (bool_para, float_para, str_para, int_para, list_para) so they land at
# the same shared-string table indices as the target workbook.
$ws.Range("E2").Value = "bool_para"
$ws.Range("E4").Value = "float_para"
$ws.Range("E6").Value = "str_para"
$ws.Range("E5").Value = "int_para"
$ws.Range("E16").Value = "list_para"

# Remaining rows reuse the now-existing shared strings.
$ws.Range("E3").Value = "bool_para"
$ws.Range("E12").Value = "float_para"
$ws.Range("E13").Value = "float_para"
$ws.Range("E14").Value = "float_para"
$ws.Range("E15").Value = "int_para"
$ws.Range("E25").Value = "bool_para"
$ws.Range("E26").Value = "int_para"
$ws.Range("E27").Value = "float_para"
$ws.Range("E28").Value = "float_para"
$ws.Range("E29").Value = "float_para"

# Restore the previous top-pane selection, then select the new bottom-pane
# active cell last so it becomes the recorded active selection.
$ws.Range("B1").Select()
$ws.Range("H25").Select()
